# Blake_Jarwin_TE_2018.xlsx - "ran scraper to update the data"
#
# Inserts two new stat columns ("height", "weight") between the existing
# "fumbles" column (D) and "fantasy points" column (old E). The existing
# "fantasy points" column and its data shift right to column G; the two
# new columns are filled with constant scraped values (same value on
# every player row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "fantasy points" column (E) two to the right (-> G), opening
# up blank columns E:F for the new "height"/"weight" fields. This carries
# the header cell's style (bold/centered/bordered) along with it, same as
# the data that follows.
$ws.Columns("E:F").Insert()

# New headers (shared strings "height" / "weight"). The Insert() above
# already carried the old E1 header's style (bold/centered/bordered) onto
# the new E1:F1 cells, matching the rest of row 1.
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"

# Scraped height/weight are constant across every row in this sheet.
$ws.Range("E2:E17").Value = 6.416666666666667
$ws.Range("F2:F17").Value = 260
